$d = $word.ActiveDocument

# --- Change 1: "Am Anfang" -> "Zu Anfang" ---
$d.Content.Find.Execute(
    "Am Anfang", $true, $false, $false, $false, $false,
    $true, 1, $false, "Zu Anfang", 2
) | Out-Null

# --- Change 2: append " für die Befreiung und Aufbau der DDR" before the
#     period that ends the "...dankt auch der Sowjetunion." sentence ---
$d.Content.Find.Execute(
    "dankt auch der Sowjetunion.", $true, $false, $false, $false, $false,
    $true, 1, $false, "dankt auch der Sowjetunion für die Befreiung und Aufbau der DDR.", 2
) | Out-Null

# --- Change 3: append four new paragraphs at the end of the document ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.Text = "!!!Mehr ins Detail gehen"

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
# p2 stays empty

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.Text = "Intention:"

$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$p4.Range.Text = "Mit dieser Rede versucht Erich Honecker in der DDR, ein Land das zu diesem Zeitpunkt in Chaos liegt wieder Vertrauen in die Regierung und in das kommunistische System zu erwecken indem er vor einem ausgewählten Publikum über die verschiedenen Erfolge und Errungenschaften spricht sowie mistrauen über die Westlichen Länder, ins besonders der BRD zu verbreiten."
